# Atualização de bases das ligas, do dia: 12-06-2024 às 23:38
# Rows 117-120 (sheet rows for ids 115-118) had their data rows cyclically
# rotated: the old row 120 data moved up to row 117, and rows 117-119 each
# shifted down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B..AD (2..30), in order, for each destination row.
# Row 117 <- old row 120 ; Row 118 <- old row 117 ; Row 119 <- old row 118 ; Row 120 <- old row 119
$rows = @{
    117 = @(7013885, "Uruguay Primera División", 45267.70833333334, "La Luz", "Atletico Fenix Montevideo", 0, 2, 0, 0, "A", 3, 3, 2.4, 2.9, 2.75, 2.6, 0, 2.025, 1.825, 2, 2.025, 1.825, -1, -1, 1.6, -1, 0.825, 0, 0)
    118 = @(7013409, "Uruguay Primera División", 45267.70833333334, "Nacional De Football", "Torque", 1, 1, 0, 1, "D", 1.666, 3.9, 4.5, 1.615, 4, 4.75, -0.75, 1.8, 2.05, 2.75, 1.95, 1.9, -1, 3, -1, -1, 1.05, -1, 0.8999999999999999)
    119 = @(7013702, "Uruguay Primera División", 45267.70833333334, "Defensor Sporting", "Danubio", 0, 2, 0, 1, "A", 1.8, 3.6, 4.2, 1.8, 3.6, 4.2, -0.75, 2.05, 1.8, 2.25, 1.85, 2, -1, -1, 3.2, -1, 0.8, -0.5, 0.5)
    120 = @(7013886, "Uruguay Primera División", 45267.70833333334, "Racing Club de Montevideo", "Cerro", 0, 1, 0, 0, "A", 2.25, 3.1, 3.25, 2.25, 2.875, 3.5, -0.25, 1.95, 1.9, 2, 1.925, 1.925, -1, -1, 2.5, -1, 0.8999999999999999, -1, 0.925)
}

foreach ($r in $rows.Keys) {
    $values = $rows[$r]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 2   # start at column B (2)
        $ws.Cells.Item($r, $col).Value = $values[$i]
    }
}
